$wb = $excel.ActiveWorkbook

# Sheet "Estadisticos 1P" - update row 5 (group 1FM)
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Cells.Item(5, 4).Value = 2
$ws1.Cells.Item(5, 6).Value = 22
$ws1.Cells.Item(5, 7).Value = 91.67

# Sheet "Estadisticos 2P" - update row 5 (group 1FM)
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Cells.Item(5, 5).Value = 22

# Sheet "Estadisticos Final" - update row 5 (group 1FM)
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Cells.Item(5, 4).Value = 2
$ws3.Cells.Item(5, 6).Value = 22
$ws3.Cells.Item(5, 7).Value = 91.67

# Sheet "Rescatables" - add new rows
$ws4 = $wb.Worksheets.Item("Rescatables")

$colA = @(21330051920277, 21330051920284, 21330051920288, 21330051920306, 21330051920324, 21330051920225, 21330051920271, 21330051920305)
$colB = @("ROMERO", "VASQUEZ", "VELAZQUEZ", "LOPEZ", "ROJAS", "XOTLANIHUA", "PEREZ", "JUSTO")
$colC = @("ZEPEDA", "ARELLANO", "TEXCAHUA", "RAMOS", "VENEGAS", "ESPINOSA", "APONTE", "NEGRETE")
$colD = @("BRYAN ABRAHAM", "RICARDO", "NANCY PAOLA", "ANETTE JOCELYN", "ANDRIK YOSIMAR", "MIXTLI TONATI", "NATALIA", "MARIA MICHELLE")
$colE = @("ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA", "ÁLGEBRA")
$colF = @("1DM", "1DM", "1DM", "1EM", "1EM", "1CM", "1DM", "1EM")
$colG = @(6, 6, 6, 6, 6, 6, 6, 6)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws4.Cells.Item($i + 2, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws4.Cells.Item($i + 2, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws4.Cells.Item($i + 2, 5).Value = $colE[$i]
}
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws4.Cells.Item($i + 2, 6).Value = $colF[$i]
}
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws4.Cells.Item($i + 2, 7).Value = $colG[$i]
}
